# Apply cryptos list price/volume updates (and one row swap at 49/50)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'34.976.07"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").Value = "'1.845.87"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.08%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'232.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.78%  "
$ws.Range("D6").Value = "'0.619"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.16%  "
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("D8").Value = "'41.84"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +6.72%  "
$ws.Range("D9").Value = "'0.328"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.27%  "
$ws.Range("D10").Value = "'0.0694"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.30%  "
$ws.Range("D11").Value = "'0.0983"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.92%  "
$ws.Range("D12").Value = "'2.113.56"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.10%  "
$ws.Range("D13").Value = "'11.44"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.33%  "
$ws.Range("D14").Value = "'1.841.59"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.91%  "
$ws.Range("E15").Value = "  +2.08%  "
$ws.Range("E16").Value = "  +2.91%  "
$ws.Range("D17").Value = "'34.982.02"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.40%  "
$ws.Range("D18").Value = "'70.06"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.16%  "
$ws.Range("E19").Value = "  +1.61%  "
$ws.Range("D20").Value = "'240.76"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.70%  "
$ws.Range("E21").Value = "  +4.10%  "
$ws.Range("E22").Value = "  +2.82%  "
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("E24").Value = "  +0.85%  "
$ws.Range("D25").Value = "'172.48"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.57%  "
$ws.Range("E26").Value = "  +0.71%  "
$ws.Range("D27").Value = "'17.51"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.92%  "
$ws.Range("D28").Value = "'0.124"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.96%  "
$ws.Range("D29").Value = "'1.64"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +8.41%  "
$ws.Range("E30").Value = "  +0.14%  "
$ws.Range("D31").Value = "'0.0553"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.28%  "
$ws.Range("E32").Value = "  -0.39%  "
$ws.Range("E33").Value = "  -0.13%  "
$ws.Range("D34").Value = "'1.64"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +22.18%  "
$ws.Range("E35").Value = "  +11.11%  "
$ws.Range("D36").Value = "'0.749"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +9.96%  "
$ws.Range("D37").Value = "'1.23"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.11%  "
$ws.Range("E38").Value = "  +11.92%  "
$ws.Range("D39").Value = "'89.88"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.86%  "
$ws.Range("D40").Value = "'1.349.86"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.11%  "
$ws.Range("E41").Value = "  +2.87%  "
$ws.Range("D42").Value = "'14.61"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.49%  "
$ws.Range("D43").Value = "'2.29"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.09%  "
$ws.Range("E44").Value = "  -1.48%  "
$ws.Range("E45").Value = "  +2.30%  "
$ws.Range("E46").Value = "  +4.40%  "
$ws.Range("D47").Value = "'6.34"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.67%  "
$ws.Range("D48").Value = "'2.032.45"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.04%  "
$ws.Range("B49").Value = "PaxDollar"
$ws.Range("C49").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D49").Value = "'1.01"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.20%  "
$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D50").Value = "'3.40"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +16.12%  "
$ws.Range("E51").Value = "  -0.20%  "
